$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank row at row 3 - shifts everything from old row 3 onward down by one.
$ws.Rows(3).Insert()

# 2) New label next to the "Ergebnisse: " header (G2) -> H2
$ws.Range("H2").Value = "Durchschnitt der F1-Werte"

# 3) New sub-heading in the now-empty row 3
$ws.Range("B3").Value = "Netz: InceptionV3"
$ws.Range("B3").Font.Bold = $false
